{"js": "// Translate the resume title and the three job-date ranges from\n// English/mixed formatting into German formatting.\nconst replacements = [\n  [\"Resume - Patti Fernandez\", \"Lebenslauf: Patti Fernandez\"],\n  [\n    \"ABC Studios: Lead Animator (Jan 2018 - Present)\",\n    \"ABC Studios: Lead Animator (Jan. 2018 - heute)\",\n  ],\n  [\n    \"XYZ Media: Senior Animator (Jun 2015 - Dez 2017)\",\n    \"XYZ Media: Senior Animator (Jun. 2015 - Dez. 2017)\",\n  ],\n  [\n    \"MNO Entertainment: Junior Animator (Sep 2012 - Mai 2015)\",\n    \"MNO Entertainment: Junior Animator (Sept. 2012 - Mai 2015)\",\n  ],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Translate the resume title and the three job-date ranges from\n# English/mixed formatting into German formatting.\n$d = $word.ActiveDocument\n\nfunction Replace-AllText {\n    param(\n        [string]$FindText,\n        [string]$ReplaceText\n    )\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute(\n        $FindText,      # FindText\n        $true,          # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap = wdFindContinue\n        $false,         # Format\n        $ReplaceText,   # ReplaceWith\n        2               # Replace = wdReplaceAll\n    )\n\n    if (-not $found) {\n        throw \"Text not found: $FindText\"\n    }\n}\n\nReplace-AllText \"Resume - Patti Fernandez\" \"Lebenslauf: Patti Fernandez\"\nReplace-AllText \"ABC Studios: Lead Animator (Jan 2018 - Present)\" \"ABC Studios: Lead Animator (Jan. 2018 - heute)\"\nReplace-AllText \"XYZ Media: Senior Animator (Jun 2015 - Dez 2017)\" \"XYZ Media: Senior Animator (Jun. 2015 - Dez. 2017)\"\nReplace-AllText \"MNO Entertainment: Junior Animator (Sep 2012 - Mai 2015)\" \"MNO Entertainment: Junior Animator (Sept. 2012 - Mai 2015)\"\n"}
